# 62. Build Util functions to read and update excel file strategically complete
# Data fix-up on the "fruit_name" column:
#   row 3 (sno=2): "Apple"  -> "iPhone"
#   row 5 (sno=4): "Banana" -> "Republic"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "iPhone"
$ws.Range("B5").Value = "Republic"

$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("B4").Select()
